$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23+ down by one.
$ws.Rows("23:23").Insert()

# Populate the new row 23 with the Salmon Recovery Conference data.
$ws.Range("A23").Value = "Salmon Recoverry Conference"
$ws.Range("B23").Value = "N/A"
$ws.Range("C23").Value = 43198
$ws.Range("D23").Value = 43199
$ws.Range("E23").Value = 2019
$ws.Range("F23").Value = "Tacoma"
$ws.Range("G23").Value = "Washington"
$ws.Range("H23").Value = "Water"
$ws.Range("I23").Value = "Puget Sound"
$ws.Range("J23").Value = 800
$ws.Range("N23").Value = "https://www.rco.wa.gov/salmon_recovery/2019-SalmonConference/Confhome.shtml"

# Move the selection to N23, matching the author's final cursor position.
$ws.Range("N23").Select() | Out-Null
